$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '29.139.38'
$ws.Range('E2').Value = '  +0.01%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.832.67'
$ws.Range('E3').Value = '  -0.33%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9996'
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '241.38'
$ws.Range('E5').Value = '  +0.64%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.6633'
$ws.Range('E6').Value = '  -2.64%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.07425'
$ws.Range('E8').Value = '  -0.35%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.2937'
$ws.Range('E9').Value = '  -1.78%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '22.71'
$ws.Range('E10').Value = '  -2.19%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07739'
$ws.Range('E11').Value = '  +1.24%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.844.53'
$ws.Range('E12').Value = '  +0.27%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.989'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.6688'
$ws.Range('E14').Value = '  -1.74%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '82.81'
$ws.Range('E15').Value = '  -5.28%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '6.100'
$ws.Range('E16').Value = '  -0.86%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.000008363'
$ws.Range('E17').Value = '  +1.85%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '29.164.91'
$ws.Range('E18').Value = '  +0.07%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '227.29'
$ws.Range('E19').Value = '  -1.39%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.48'
$ws.Range('E20').Value = '  -0.12%  '
$ws.Range('E21').Value = '  +0.12%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '7.158'
$ws.Range('E22').Value = '  -2.56%  '
$ws.Range('E23').Value = '  +0.04%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '159.77'
$ws.Range('E24').Value = '  -0.86%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '8.622'
$ws.Range('E25').Value = '  -0.96%  '
$ws.Range('E26').Value = '  -1.78%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '17.97'
$ws.Range('E27').Value = '  -0.53%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.508'
$ws.Range('E28').Value = '  +0.33%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '4.115'
$ws.Range('E29').Value = '  -3.24%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.043'
$ws.Range('E30').Value = '  -2.23%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.194'
$ws.Range('E31').Value = '  -0.09%  '
$ws.Range('E32').Value = '  -0.44%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.870'
$ws.Range('E33').Value = '  +1.24%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.7524'
$ws.Range('E34').Value = '  -0.20%  '
$ws.Range('E35').Value = '  +0.39%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.630'
$ws.Range('E36').Value = '  -1.97%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.282.75'
$ws.Range('E37').Value = '  -2.22%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01797'
$ws.Range('E38').Value = '  -1.60%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.733'
$ws.Range('E39').Value = '  +0.63%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.9288'
$ws.Range('E40').Value = '  -1.55%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.08783'
$ws.Range('E41').Value = '  +10.22%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.941'
$ws.Range('E42').Value = '  -2.06%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.000'
$ws.Range('E43').Value = '  +0.11%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '102.07'
$ws.Range('E44').Value = '  -2.63%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.977.82'
$ws.Range('E45').Value = '  -0.49%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.5149'
$ws.Range('E46').Value = '  -0.59%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.767'
$ws.Range('E47').Value = '  -0.35%  '
$ws.Range('E48').Value = '  -1.14%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '63.29'
$ws.Range('E49').Value = '  -1.10%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.05900'
$ws.Range('E50').Value = '  -0.83%  '
$ws.Range('B51').Value = 'EnergySwap'
$ws.Range('C51').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '8.830'
$ws.Range('E51').Value = '  -6.35%  '
